$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.284680485725403
$ws.Range("B1").Value = 2.31641149520874
$ws.Range("C1").Value = 3.956624984741211
$ws.Range("D1").Value = 2.994519710540771
$ws.Range("E1").Value = 1.335245251655579
